# Apply the "Saldo" export update:
#  - add a new account row (VITOR) before THOMAS' row
#  - remove WASHINGTON's row
#  - add a new account row (MARCIA) before CARLOS' row
#  - add a new account row (BLUEMETRIX) before MARCO's row
#  - add a new account row (NILBORN) right after VIVIANE's row
#
# Helpers insert a brand-new data row directly above/below the row
# currently holding a given account number in column A, then fill it in.
# Column A is forced to text ("@") so the zero-padded account number
# keeps its leading zeros instead of being coerced to a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Insert-AcctRowBefore {
    param([string]$anchorAccount, [string]$account, [string]$name, [double]$saldo)

    $anchorCell = $ws.Columns("A:A").Find($anchorAccount)
    $targetRow = $anchorCell.Row

    $ws.Rows("$targetRow`:$targetRow").Insert()

    # Only column A needs to be pinned to Text - it's the only column whose
    # value looks like a number (zero-padded account id); the name column
    # is already plain text and the balance column is meant to be numeric.
    $ws.Range("A$targetRow").NumberFormat = "@"
    $ws.Range("A$targetRow").Value = $account
    $ws.Range("B$targetRow").Value = $name
    $ws.Range("C$targetRow").Value = $saldo
}

function Insert-AcctRowAfter {
    param([string]$anchorAccount, [string]$account, [string]$name, [double]$saldo)

    $anchorCell = $ws.Columns("A:A").Find($anchorAccount)
    $targetRow = $anchorCell.Row + 1

    $ws.Rows("$targetRow`:$targetRow").Insert()

    $ws.Range("A$targetRow").NumberFormat = "@"
    $ws.Range("A$targetRow").Value = $account
    $ws.Range("B$targetRow").Value = $name
    $ws.Range("C$targetRow").Value = $saldo
}

function Remove-AcctRow {
    param([string]$account)

    $cell = $ws.Columns("A:A").Find($account)
    $row = $cell.Row
    $ws.Rows("$row`:$row").Delete()
}

# 1. New row for VITOR, right before THOMAS (004224011)
Insert-AcctRowBefore "004224011" "002694089" "VITOR" 45548.27

# 2. Remove WASHINGTON's row
Remove-AcctRow "005231126"

# 3. New row for MARCIA, right before CARLOS (004488571)
Insert-AcctRowBefore "004488571" "005203796" "MARCIA" 3739.5

# 4. New row for BLUEMETRIX, right before MARCO (004435987)
Insert-AcctRowBefore "004435987" "001761119" "BLUEMETRIX" 126.98

# 5. New row for NILBORN, right after VIVIANE (001294033 / 79.82)
Insert-AcctRowAfter "001294033" "005073033" "NILBORN" 79.39
